# edit.ps1 - COM-interop script reproducing commit "Update countries & provincias Spain"
#
# The source refresh re-sorted a handful of countries that were tied (or nearly
# tied) on "Casos totales", which shifts their row position by one or two slots
# in a few local clusters; every other row keeps its numbers but gets the fresh
# totals pulled from the new snapshot. We reproduce both effects by writing the
# final (country name + B:H numbers) directly into each affected row, plus the
# "updated as of" footer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer: "Datos actualizados a ..." timestamp lives in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Septiembre de 2020 a las 17:08"

# Row 4
$ws.Range("B4").Value = 7250345
$ws.Range("C4").Value = 6161
$ws.Range("D4").Value = 4481095
$ws.Range("E4").Value = 2560640
$ws.Range("G4").Value = 170
$ws.Range("H4").Value = 208610

# Row 5
$ws.Range("B5").Value = 5915753
$ws.Range("C5").Value = 14182
$ws.Range("D5").Value = 4852313
$ws.Range("E5").Value = 969979
$ws.Range("G5").Value = 51
$ws.Range("H5").Value = 93461

# Row 15
$ws.Range("B15").Value = 455979
$ws.Range("C15").Value = 2111
$ws.Range("D15").Value = 430259
$ws.Range("E15").Value = 13129
$ws.Range("G15").Value = 64
$ws.Range("H15").Value = 12591

# Row 23
$ws.Range("B23").Value = 308104
$ws.Range("C23").Value = 1869
$ws.Range("D23").Value = 223693
$ws.Range("E23").Value = 48593
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = 35818

# Row 29
$ws.Range("B29").Value = 150891
$ws.Range("C29").Value = 435
$ws.Range("D29").Value = 129911
$ws.Range("E29").Value = 11724
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 9256

# Row 51  (A51 -> Portugal)
$ws.Range("A51").Value = "Portugal"
$ws.Range("B51").Value = 72939
$ws.Range("C51").Value = 884
$ws.Range("D51").Value = 47380
$ws.Range("E51").Value = 23615
$ws.Range("G51").Value = 8
$ws.Range("H51").Value = 1944

# Row 52  (A52 -> Etiopia)
$ws.Range("A52").Value = "Etiopia"
$ws.Range("B52").Value = 72173
$ws.Range("D52").Value = 29863
$ws.Range("E52").Value = 41155
$ws.Range("H52").Value = 1155

# Row 59
$ws.Range("D59").Value = 57359
$ws.Range("E59").Value = 299

# Row 86
$ws.Range("B86").Value = 17483
$ws.Range("C86").Value = 140
$ws.Range("D86").Value = 14516
$ws.Range("E86").Value = 2245
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 722

# Row 95
$ws.Range("B95").Value = 13153
$ws.Range("C95").Value = 108
$ws.Range("D95").Value = 7397
$ws.Range("E95").Value = 5381
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 375

# Row 96
$ws.Range("B96").Value = 10918
$ws.Range("C96").Value = 83
$ws.Range("D96").Value = 8749
$ws.Range("E96").Value = 2049
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 120

# Row 101  (A101 -> Birmania)
$ws.Range("A101").Value = "Birmania"
$ws.Range("B101").Value = 9991
$ws.Range("C101").Value = 879
$ws.Range("D101").Value = 2681
$ws.Range("E101").Value = 7112
$ws.Range("G101").Value = 24
$ws.Range("H101").Value = 198

# Row 102  (A102 -> Montenegro)
$ws.Range("A102").Value = "Montenegro"
$ws.Range("B102").Value = 9962
$ws.Range("D102").Value = 6177
$ws.Range("E102").Value = 3630
$ws.Range("H102").Value = 155

# Row 103  (A103 -> Guayana Francesa)
$ws.Range("A103").Value = "Guayana Francesa"
$ws.Range("B103").Value = 9831
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 9472
$ws.Range("E103").Value = 294
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 65

# Row 104  (A104 -> Tayikistan)
$ws.Range("A104").Value = "Tayikistan"
$ws.Range("B104").Value = 9605
$ws.Range("C104").Value = 43
$ws.Range("D104").Value = 8385
$ws.Range("E104").Value = 1145
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 75

# Row 105  (A105 -> Finlandia)
$ws.Range("A105").Value = "Finlandia"
$ws.Range("B105").Value = 9577
$ws.Range("C105").Value = 93
$ws.Range("D105").Value = 7850
$ws.Range("E105").Value = 1384
$ws.Range("H105").Value = 343

# Row 107
$ws.Range("B107").Value = 8723
$ws.Range("C107").Value = 39
$ws.Range("E107").Value = 1945

# Row 110  (A110 -> Jordania)
$ws.Range("A110").Value = "Jordania"
$ws.Range("B110").Value = 8061
$ws.Range("C110").Value = 850
$ws.Range("D110").Value = 4131
$ws.Range("E110").Value = 3887
$ws.Range("G110").Value = 4
$ws.Range("H110").Value = 43

# Row 111  (A111 -> Zimbabue)
$ws.Range("A111").Value = "Zimbabue"
$ws.Range("B111").Value = 7787
$ws.Range("D111").Value = 6057
$ws.Range("E111").Value = 1503
$ws.Range("H111").Value = 227

# Row 112  (A112 -> Mozambique)
$ws.Range("A112").Value = "Mozambique"
$ws.Range("B112").Value = 7589
$ws.Range("D112").Value = 4649
$ws.Range("E112").Value = 2887
$ws.Range("H112").Value = 53

# Row 113  (A113 -> Mauritania)
$ws.Range("A113").Value = "Mauritania"
$ws.Range("B113").Value = 7457
$ws.Range("D113").Value = 7070
$ws.Range("E113").Value = 226
$ws.Range("H113").Value = 161

# Row 114  (A114 -> Uganda)
$ws.Range("A114").Value = "Uganda"
$ws.Range("B114").Value = 7364
$ws.Range("C114").Value = 146
$ws.Range("D114").Value = 3647
$ws.Range("E114").Value = 3646
$ws.Range("H114").Value = 71

# Row 115  (A115 -> Jamaica)
$ws.Range("A115").Value = "Jamaica"
$ws.Range("B115").Value = 5854
$ws.Range("C115").Value = 131
$ws.Range("D115").Value = 1624
$ws.Range("E115").Value = 4142
$ws.Range("G115").Value = 8
$ws.Range("H115").Value = 88

# Row 116  (A116 -> Malaui)
$ws.Range("A116").Value = "Malaui"
$ws.Range("B116").Value = 5764
$ws.Range("D116").Value = 4178
$ws.Range("E116").Value = 1407
$ws.Range("H116").Value = 179

# Row 133
$ws.Range("B133").Value = 4285
$ws.Range("C133").Value = 8
$ws.Range("D133").Value = 2185
$ws.Range("E133").Value = 2030

# Row 206  (A206 -> Timor Oriental)
$ws.Range("A206").Value = "Timor Oriental"

# Row 207  (A207 -> Santa Lucia)
$ws.Range("A207").Value = "Santa Lucia"
